# Conclusion.docx edit script
# Applies the changes described by the commit: reworks the closing
# paragraphs of the thesis conclusion, adding a new paragraph about the
# target audience and tightening the wording of the final paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph "ALAT is to be the new recommended authoring environment
#    for GALE. ...": tweak the final sentence about templating/knowledge.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "usability and to author adaptivity without knowledge of",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "usability. It also allows authors to author adaptivity without any required knowledge on",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) Insert a brand-new paragraph about the target audience right after
#    the paragraph that now ends in "... adaptation code."
# ---------------------------------------------------------------------
$targetAudienceParagraph = $d.Paragraphs.Item(3)
$targetAudienceParagraph.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(4)
$newPara.Range.Text = "The target audience consists out of students and engineers with at least a basic understanding of adaptive hypermedia. A result of this is that ALAT is more difficult to author with for non-academic users. This is due to the complexity and user experience which is a result of ALAT" + [char]0x2019 + "s extensive generic behavior. A version of ALAT with a limited scope, stripped of some more advanced features would lower this barrier-of-entry and could make ALAT more suitable for non-experts and companies such as " + [char]0x201C + "De Roode Kikker" + [char]0x201D + "."

# ---------------------------------------------------------------------
# 3) Rework the final paragraph: swap the opening, fix a preposition,
#    change "serves as" -> "provides", add a new sentence about target
#    audience/genericity, and tighten the closing sentence's wording.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "This new authoring tool contributes to the usability of GALE and brings a new player in the field of adaptive hypermedia authoring. This thesis serves as",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ALAT contributes to the usability of GALE and brings a new player to the field of adaptive hypermedia authoring. This thesis provides",
    2) | Out-Null

$d.Content.Find.Execute(
    "adaptive hypermedia authoring. ALAT explores authoring in academic use and innovates by combining an interface with simple controls with extensive templating.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "adaptive hypermedia authoring. The main contributing factors to these differences are target audience and genericity. ALAT explores authoring by academic users and innovates by combining a simple user interface with extensive templating.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4) Relocate the "_GoBack" last-edit bookmark to sit right after
#    "academic users", mirroring where the author's cursor ended up.
# ---------------------------------------------------------------------
$cursorMark = $d.Content
$cursorMark.Find.Execute("by academic users", $true) | Out-Null
$cursorPos = $cursorMark.End
$d.Bookmarks.Add("_GoBack", $d.Range($cursorPos, $cursorPos))
